$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Temporarily force the Price column to Text so numeric-looking values
# (e.g. "258.12") are stored as strings, matching the source data which
# always uses inline/shared strings for these cells.
$ws.Range("D2:D50").NumberFormat = "@"

$ws.Range("D2").Value = "98.459.62"
$ws.Range("E2").Value = "  -0.15%  "
$ws.Range("D3").Value = "3.367.70"
$ws.Range("E3").Value = "  +0.18%  "
$ws.Range("E4").Value = "  -0.07%  "
$ws.Range("D5").Value = "258.12"
$ws.Range("E5").Value = "  -0.15%  "
$ws.Range("D6").Value = "666.44"
$ws.Range("E6").Value = "  +5.99%  "
$ws.Range("E7").Value = "  +10.48%  "
$ws.Range("D8").Value = "0.466"
$ws.Range("E8").Value = "  +19.70%  "
$ws.Range("E9").Value = "  +25.26%  "
$ws.Range("E10").Value = "  -0.03%  "
$ws.Range("D11").Value = "3.364.50"
$ws.Range("E11").Value = "  +0.24%  "
$ws.Range("D12").Value = "0.212"
$ws.Range("E12").Value = "  +6.78%  "
$ws.Range("D13").Value = "42.29"
$ws.Range("E13").Value = "  +16.03%  "
$ws.Range("D14").Value = "0.0000271"
$ws.Range("E14").Value = "  +9.46%  "
$ws.Range("D15").Value = "99.113.40"
$ws.Range("E15").Value = "  +0.83%  "
$ws.Range("D16").Value = "3.991.25"
$ws.Range("E16").Value = "  +0.12%  "
$ws.Range("E17").Value = "  +3.05%  "
$ws.Range("D18").Value = "3.372.10"
$ws.Range("E18").Value = "  +0.41%  "
$ws.Range("D19").Value = "7.64"
$ws.Range("E19").Value = "  +25.86%  "
$ws.Range("D20").Value = "16.91"
$ws.Range("E20").Value = "  +11.12%  "
$ws.Range("E21").Value = "  +1.24%  "
$ws.Range("D22").Value = "531.70"
$ws.Range("E22").Value = "  +8.92%  "
$ws.Range("D23").Value = "10.53"
$ws.Range("E23").Value = "  +11.85%  "
$ws.Range("D24").Value = "0.0000219"
$ws.Range("E24").Value = "  +4.67%  "
$ws.Range("D25").Value = "0.437"
$ws.Range("E25").Value = "  +54.89%  "
$ws.Range("D26").Value = "102.77"
$ws.Range("E26").Value = "  +15.37%  "
$ws.Range("D27").Value = "6.26"
$ws.Range("E27").Value = "  +10.92%  "
$ws.Range("D28").Value = "12.64"
$ws.Range("E28").Value = "  +6.15%  "
$ws.Range("D29").Value = "3.546.67"
$ws.Range("E29").Value = "  -0.06%  "
$ws.Range("E30").Value = "  +10.34%  "
$ws.Range("D31").Value = "0.999"
$ws.Range("E31").Value = "  -0.06%  "
$ws.Range("D32").Value = "11.07"
$ws.Range("E32").Value = "  +14.25%  "
$ws.Range("E33").Value = "  -0.26%  "
$ws.Range("D34").Value = "0.997"
$ws.Range("E34").Value = "  -0.51%  "
$ws.Range("D35").Value = "29.57"
$ws.Range("E35").Value = "  +5.20%  "
$ws.Range("D36").Value = "0.552"
$ws.Range("E36").Value = "  +20.02%  "
$ws.Range("D37").Value = "7.90"
$ws.Range("E37").Value = "  +8.56%  "
$ws.Range("E38").Value = "  +8.51%  "
$ws.Range("E39").Value = "  +5.29%  "
$ws.Range("D40").Value = "529.47"
$ws.Range("E40").Value = "  +6.17%  "
$ws.Range("E41").Value = "  +6.51%  "
$ws.Range("E42").Value = "  -0.46%  "
$ws.Range("E43").Value = "  +34.45%  "
$ws.Range("E44").Value = "  +1.81%  "
$ws.Range("B45").Value = "ARBITRUM"
$ws.Range("C45").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D45").Value = "0.842"
$ws.Range("E45").Value = "  +6.20%  "
$ws.Range("B46").Value = "dogwifhat"
$ws.Range("C46").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$ws.Range("D46").Value = "3.44"
$ws.Range("E46").Value = "  +4.62%  "
$ws.Range("E47").Value = "  +0.03%  "
$ws.Range("B48").Value = "Cosmos"
$ws.Range("C48").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D48").Value = "7.99"
$ws.Range("E48").Value = "  +19.53%  "
$ws.Range("B49").Value = "Stacks"
$ws.Range("C49").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D49").Value = "2.08"
$ws.Range("E49").Value = "  +7.31%  "
$ws.Range("D50").Value = "5.14"
$ws.Range("E50").Value = "  +11.43%  "
$ws.Range("E51").Value = "  +11.71%  "

# Restore default (no explicit number format) so styling matches original.
$ws.Range("D2:D50").ClearFormats()
